$d = $word.ActiveDocument

# Hybrid bold + color (2C3E50) highlighting for quantitative impact metrics.
# Applied to specific paragraphs (by 1-based Paragraphs index) so that
# look-alike numbers elsewhere in the resume (summary / key-projects
# sections) are left untouched, matching the target diff precisely.

$plusMinus = [char]0x00B1

function Highlight-Metric($paraIndex, $text) {
    $p = $d.Paragraphs.Item($paraIndex)
    $rng = $p.Range
    $found = $rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.Bold = 1
        $rng.Font.Color = 5258796
    }
}

# Para 9 (Siege Analytics bullet 1): "... from 23% to 64%"
Highlight-Metric 9 "23%"
Highlight-Metric 9 "64%"

# Para 11 (Siege Analytics bullet 3): "Achieved 87% ... 71%, ... from ±4.2% to ±2.1%"
Highlight-Metric 11 "87%"
Highlight-Metric 11 "71%"
$negFourTwo = $plusMinus + "4.2%"
$negTwoOne = $plusMinus + "2.1%"
Highlight-Metric 11 $negFourTwo
Highlight-Metric 11 $negTwoOne

# Para 31 (Myers Research bullet 3): "... bids from 1,200 vendors ..."
Highlight-Metric 31 "1,200"

# Para 46 (Lake Research Partners bullet 3): "... $400M ... $1B+"
Highlight-Metric 46 '$400M'
Highlight-Metric 46 '$1B'

# Para 63 (Key Achievements bullet 2): "... by 73.5%, ... $4.7M"
Highlight-Metric 63 "73.5%"
Highlight-Metric 63 '$4.7M'

# Para 65 (Key Achievements bullet 4): "Achieved 87% ... of 71%"
Highlight-Metric 65 "87%"
Highlight-Metric 65 "71%"

Write-Output "done"
